$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab (and underlying workbook sheet name)
$ws.Name = "causa_muerte"

# Remove the header styling (bold white font on blue fill, centered) by
# clearing formatting on the header row so cells fall back to the default style
$headerRange = $ws.Range("A1:D1")
$headerRange.ClearFormats()

# Update header text to the new lowercase/underscored values
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "tipo_causa"
$ws.Range("D1").Value = "comentario"
